$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new blank (formatted) row is inserted above the "Indice gulpease" row,
# which pushes it (and the trailing blank formatting row) one row down.
$ws.Rows("11:11").Insert()

# Fill in the previously empty "AR" (column G) verification values.
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 5
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 0

# Add the new "total" column (I) with per-row sums.
$ws.Range("I2").Formula = "=SUM(B2+E2+F2+G2,C2+D2)"
$ws.Range("I3").Formula = "=SUM(B3+C3+D3+E3+F3+G3)"
$ws.Range("I4").Formula = "=SUM(B4+C4+D4+E4+F4+G4)"
$ws.Range("I5").Formula = "=SUM(B5+C5+D5+E5+F5+G5)"
$ws.Range("I6").Formula = "=SUM(B6+C6+D6+E6+F6+G6)"
$ws.Range("I7").Value = 9
$ws.Range("I8").Value = 9
$ws.Range("I9").Formula = "=SUM(B9+C9+D9+E9+F9+G9)"
$ws.Range("I10").Formula = "=SUM(B10+C10+D10+E10+F10+G10)"

$ws.Range("I10").Select() | Out-Null
